$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "last updated" timestamp footer (row 1) ---
$ws.Range("A1").Value = "Datos actualizados a 27 de Abril de 2020 a las 21:22"

# --- Row 4: Estados Unidos - refreshed totals ---
$ws.Range("B4").Value = 1000441
$ws.Range("C4").Value = 13281
$ws.Range("D4").Value = 137465
$ws.Range("E4").Value = 806690
$ws.Range("F4").Value = 14175
$ws.Range("G4").Value = 873
$ws.Range("H4").Value = 56286

# --- Rows 99/100: Costa Rica overtakes Niger in the ranking ---
$ws.Range("A99").Value = "Costa Rica"
$ws.Range("B99").Value = 697
$ws.Range("C99").Value = 2
$ws.Range("D99").Value = 287
$ws.Range("E99").Value = 404
$ws.Range("F99").Value = 8
$ws.Range("G99").Value = 0
$ws.Range("H99").Value = 6

$ws.Range("A100").Value = "Niger"
$ws.Range("B100").Value = 696
$ws.Range("C100").Value = 0
$ws.Range("D100").Value = 350
$ws.Range("E100").Value = 317
$ws.Range("F100").Value = 0
$ws.Range("G100").Value = 0
$ws.Range("H100").Value = 29

# --- Rows 131/132: Ruanda overtakes Congo in the ranking ---
$ws.Range("A131").Value = "Ruanda"
$ws.Range("B131").Value = 207
$ws.Range("C131").Value = 16
$ws.Range("D131").Value = 93
$ws.Range("E131").Value = 114
$ws.Range("F131").Value = 0
$ws.Range("G131").Value = 0
$ws.Range("H131").Value = 0

$ws.Range("A132").Value = "Congo"
$ws.Range("B132").Value = 200
$ws.Range("C132").Value = 0
$ws.Range("D132").Value = 19
$ws.Range("E132").Value = 175
$ws.Range("F132").Value = 0
$ws.Range("G132").Value = 0
$ws.Range("H132").Value = 6

# --- Row 178: Angola - refreshed totals ---
$ws.Range("B178").Value = 27
$ws.Range("C178").Value = 1
$ws.Range("D178").Value = 6
$ws.Range("E178").Value = 19
$ws.Range("F178").Value = 0
$ws.Range("G178").Value = 0
$ws.Range("H178").Value = 2
